$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header D1: bold + centered (create this style first so it lands at cellXfs index 3) ---
$ws.Range("D1").HorizontalAlignment = -4108

# --- Existing Done column cells D2:D4 - center align (plain+center -> cellXfs index 4) ---
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D4").HorizontalAlignment = -4108

# Insert a new row at position 9 ("Reach out to business"); shifts existing rows 9-17 down to 10-18
$ws.Rows("9:9").Insert()

# --- New row 9: Reach out to business ---
$ws.Range("A9").Value = "Reach out to business"
$ws.Range("B9").Value = "Angel Todev"
$ws.Range("C9").NumberFormat = $ws.Range("C2").NumberFormat()
$ws.Range("C9").Value = 43393
$ws.Range("D9").Value = "x"
$ws.Range("D9").HorizontalAlignment = -4108

# --- Row 8: Complete login system - add completion date + done mark ---
$ws.Range("C8").NumberFormat = $ws.Range("C2").NumberFormat()
$ws.Range("C8").Value = 43394
$ws.Range("D8").Value = "x"
$ws.Range("D8").HorizontalAlignment = -4108

# --- Row 10 (was row 9 pre-insert): Schedule initial meeting with business ---
$ws.Range("C10").NumberFormat = $ws.Range("C2").NumberFormat()
$ws.Range("C10").Value = 43396
$ws.Range("D10").Value = "x"
$ws.Range("D10").HorizontalAlignment = -4108

# --- Row 12 (was row 11 pre-insert): Verify that publishing methods of websites are acceptable ---
$ws.Range("C12").NumberFormat = $ws.Range("C2").NumberFormat()
$ws.Range("C12").Value = 43399
$ws.Range("D12").Value = "x"
$ws.Range("D12").HorizontalAlignment = -4108

# Match the author's final cursor position noted in the saved file
$null = $ws.Range("C21").Select()

Write-Host "done"
